# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the item has
# moved from "In Translation" to "Ready for handoff", and refreshes the
# handoff timestamps accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status + handoff date (row 2) ---
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-29-20 20:29:46"

# --- zh-cn sheet: status + handoff datetime (row 2) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-20 20:29:43"

# --- de-de sheet: status + handoff datetime (row 2) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-20 20:29:46"
